# Apply the 11-May-2023 GitHub Actions "Updated cryptos list" refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper-free, explicit writes below. Price (col D) values are written as
# TEXT (matching the source data, which stores "1.003" etc. as strings, not
# numbers) by forcing the NumberFormat to Text before the write and then
# restoring the Normal style so no stray formatting is introduced.
function Set-PriceText($addr, $text) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = "Normal"
}

# --- Rows whose Coin/Link are unchanged: refresh Price (D) and Volume(1h) (E) ---
Set-PriceText "D2" "27.231.94"
$ws.Range("E2").Value = "  -3.44%  "
Set-PriceText "D3" "1.808.01"
$ws.Range("E3").Value = "  -3.76%  "
Set-PriceText "D4" "1.003"
$ws.Range("E4").Value = "  -0.03%  "
Set-PriceText "D5" "310.80"
$ws.Range("E5").Value = "  -1.70%  "
Set-PriceText "D6" "1.002"
$ws.Range("E6").Value = "  -0.03%  "
Set-PriceText "D7" "0.4214"
$ws.Range("E7").Value = "  -2.34%  "
Set-PriceText "D8" "0.3551"
$ws.Range("E8").Value = "  -3.93%  "
Set-PriceText "D9" "0.07127"
$ws.Range("E9").Value = "  -3.98%  "
Set-PriceText "D10" "0.8484"
$ws.Range("E10").Value = "  -4.25%  "
Set-PriceText "D11" "20.18"
$ws.Range("E11").Value = "  -4.73%  "
Set-PriceText "D12" "1.801.66"
$ws.Range("E12").Value = "  -6.32%  "
Set-PriceText "D13" "5.328"
$ws.Range("E13").Value = "  -2.79%  "
Set-PriceText "D16" "1.005"
$ws.Range("E16").Value = "  -0.02%  "
Set-PriceText "D17" "81.16"
$ws.Range("E17").Value = "  -0.11%  "
Set-PriceText "D18" "0.000008754"
$ws.Range("E18").Value = "  -4.25%  "
Set-PriceText "D19" "1.004"
$ws.Range("E19").Value = "  +0.16%  "
Set-PriceText "D20" "15.10"
$ws.Range("E20").Value = "  -3.34%  "
Set-PriceText "D21" "27.666.49"
$ws.Range("E21").Value = "  -2.31%  "
Set-PriceText "D22" "5.089"
$ws.Range("E22").Value = "  -0.11%  "
Set-PriceText "D23" "10.88"
$ws.Range("E23").Value = "  -0.65%  "
Set-PriceText "D24" "2.110.21"
$ws.Range("E24").Value = "  +0.03%  "
Set-PriceText "D25" "1.962"
$ws.Range("E25").Value = "  -1.07%  "
Set-PriceText "D26" "153.37"
$ws.Range("E26").Value = "  -0.64%  "
Set-PriceText "D27" "18.24"
$ws.Range("E27").Value = "  -2.65%  "
Set-PriceText "D28" "5.063"
$ws.Range("E28").Value = "  -6.58%  "
Set-PriceText "D29" "113.06"
$ws.Range("E29").Value = "  -4.55%  "
Set-PriceText "D31" "0.08893"
$ws.Range("E31").Value = "  -1.08%  "
Set-PriceText "D35" "1.103"
$ws.Range("E35").Value = "  -6.33%  "
Set-PriceText "D36" "1.002"
$ws.Range("E36").Value = "  +0.05%  "
Set-PriceText "D38" "0.05210"
$ws.Range("E38").Value = "  -4.97%  "
Set-PriceText "D39" "0.01900"
$ws.Range("E39").Value = "  -3.63%  "
Set-PriceText "D40" "2.732"
$ws.Range("E40").Value = "  -5.77%  "
Set-PriceText "D41" "0.1638"
$ws.Range("E41").Value = "  -3.59%  "
Set-PriceText "D42" "0.4977"
$ws.Range("E42").Value = "  -3.82%  "
Set-PriceText "D43" "6.306"
$ws.Range("E43").Value = "  -8.47%  "
Set-PriceText "D44" "8.192"
$ws.Range("E44").Value = "  -4.67%  "
Set-PriceText "D47" "1.002"
$ws.Range("E47").Value = "  -0.07%  "
Set-PriceText "D48" "0.06395"
$ws.Range("E48").Value = "  -2.86%  "
Set-PriceText "D49" "0.4562"
$ws.Range("E49").Value = "  -4.42%  "
Set-PriceText "D50" "1.597"
$ws.Range("E50").Value = "  -3.81%  "
Set-PriceText "D51" "63.02"
$ws.Range("E51").Value = "  -3.72%  "

# --- Rows where only Volume(1h) (E) changed ---
$ws.Range("E30").Value = "  -9.81%  "
$ws.Range("E37").Value = "  -5.83%  "

# --- Rows that were reordered in the ranking: Coin name, Link, Price and Volume(1h) all change ---
$ws.Range("B14").Value = "Chainlink"
$ws.Range("C14").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
Set-PriceText "D14" "6.367"
$ws.Range("E14").Value = "  -3.94%  "
$ws.Range("B15").Value = "TRON"
$ws.Range("C15").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
Set-PriceText "D15" "0.06905"
$ws.Range("E15").Value = "  -1.07%  "
$ws.Range("B32").Value = "ImmutableX"
$ws.Range("C32").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
Set-PriceText "D32" "0.7421"
$ws.Range("E32").Value = "  -6.56%  "
$ws.Range("B33").Value = "Filecoin"
$ws.Range("C33").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-PriceText "D33" "4.469"
$ws.Range("E33").Value = "  -5.29%  "
$ws.Range("B34").Value = "HuobiToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
Set-PriceText "D34" "2.922"
$ws.Range("E34").Value = "  -1.17%  "
$ws.Range("B45").Value = "Quant"
$ws.Range("C45").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
Set-PriceText "D45" "105.12"
$ws.Range("E45").Value = "  -0.71%  "
$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-PriceText "D46" "10.26"
$ws.Range("E46").Value = "  -2.98%  "
